# Adds a "2022-Q4" sheet (fund-holdings detail) right after the "总计"
# summary sheet, and updates the "总计" sheet with the new quarter's
# aggregate row (shifting the existing history down by one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write a value into a cell as plain TEXT (no auto "number-ify",
# no leftover NumberFormat/style residue) so numeric-looking strings like
# fund codes ("010108") or percentages ("48.01") are preserved exactly
# as inline/shared strings, matching the source workbook's convention.
# ---------------------------------------------------------------------
function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = [string]$value
    $cell.Style = "Normal"
}

# =======================================================================
# 1) "总计" sheet: insert the new 2022-Q4 summary row, shift history down
# =======================================================================
$summary = $wb.Worksheets.Item("总计")

$summaryRows = @(
     ,@("2022-Q4", 31, 8.38)
     ,@("2022-Q3", 49, 10.36)
     ,@("2022-Q2", 18, 3.23)
     ,@("2022-Q1", 7, 1.4)
     ,@("2021-Q4", 5, 1.11)
     ,@("2021-Q3", 1, 0.57)
     ,@("2021-Q2", 6, 1.06)
     ,@("2021-Q1", 6, 1.74)
)

# Row 9 is brand new (the sheet used to stop at row 8) — clone the
# formatting of the A-column "index" cell from the row above it before
# the loop below fills in its value, so it keeps the same bold/centered/
# bordered look as A2:A8.
$summary.Cells.Item(8, 1).Copy()
$summary.Cells.Item(9, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = 2 + $i
    $vals = $summaryRows[$i]
    $summary.Cells.Item($r, 1).Value = $i
    $summary.Cells.Item($r, 2).Value = $vals[0]
    $summary.Cells.Item($r, 3).Value = $vals[1]
    $summary.Cells.Item($r, 4).Value = $vals[2]
}

# =======================================================================
# 2) New "2022-Q4" sheet: duplicate an existing fund-holdings sheet (so it
#    inherits identical formatting/styles/column layout) right after
#    "总计", then overwrite its values with the 2022-Q4 data.
# =======================================================================
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($null, $summary)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

$newSheetRows = @(
    ,@(0, "010108", "景顺长城核心招景混合A", "48.01", "89.61", "2.92", "1.4019", 9)
    ,@(1, "010027", "景顺长城核心中景一年持有期混合", "44.17", "89.32", "3.17", "1.4002", 8)
    ,@(2, "013797", "博时优质鑫选一年持有期混合A", "38.62", "79.82", "2.88", "1.1123", 9)
    ,@(3, "011756", "博时产业优选灵活配置混合A", "23.27", "84.52", "3.71", "0.8633", 4)
    ,@(4, "260116", "景顺长城核心竞争力混合A", "26.19", "85.33", "2.56", "0.6705", 9)
    ,@(5, "009190", "景顺长城核心优选一年持有期混合", "10.80", "90.64", "5.03", "0.5432", 4)
    ,@(6, "006158", "博时荣享回报灵活配置定期开放混合A", "10.78", "78.43", "3.67", "0.3956", 3)
    ,@(7, "011585", "博时产业慧选混合A", "10.03", "86.88", "3.71", "0.3721", 5)
    ,@(8, "011336", "兴全汇吉一年持有期混合A", "15.09", "39.83", "1.92", "0.2897", 5)
    ,@(9, "008866", "博时产业新趋势灵活配置混合A", "6.93", "81.37", "3.66", "0.2536", 5)
    ,@(10, "009740", "博时研究臻选三年持有期灵活配置混合A", "7.01", "80.20", "2.72", "0.1907", 8)
    ,@(11, "000936", "博时产业新动力灵活配置混合A", "4.70", "89.28", "3.81", "0.1791", 5)
    ,@(12, "009591", "博时研究精选一年持有期灵活配置混合A", "4.85", "85.76", "3.69", "0.1790", 4)
    ,@(13, "010455", "博时产业精选灵活配置混合A", "4.84", "86.08", "3.67", "0.1776", 5)
    ,@(14, "002142", "博时外延增长主题灵活配置混合", "2.28", "84.83", "3.66", "0.0834", 5)
    ,@(15, "011586", "博时产业慧选混合C", "1.01", "86.88", "3.71", "0.0375", 5)
    ,@(16, "010456", "博时产业精选灵活配置混合C", "0.72", "86.08", "3.67", "0.0264", 5)
    ,@(17, "011757", "博时产业优选灵活配置混合C", "0.70", "84.52", "3.71", "0.0260", 4)
    ,@(18, "015731", "景顺长城核心竞争力混合C", "1.01", "85.33", "2.56", "0.0259", 9)
    ,@(19, "011340", "博时战略新材料主题混合A", "0.84", "79.48", "2.53", "0.0213", 8)
    ,@(20, "009592", "博时研究精选一年持有期灵活配置混合C", "0.57", "85.76", "3.69", "0.0210", 4)
    ,@(21, "011337", "兴全汇吉一年持有期混合C", "0.80", "39.83", "1.92", "0.0154", 5)
    ,@(22, "006159", "博时荣享回报灵活配置定期开放混合C", "0.41", "78.43", "3.67", "0.0150", 3)
    ,@(23, "011341", "博时战略新材料主题混合C", "0.56", "79.48", "2.53", "0.0142", 8)
    ,@(24, "008867", "博时产业新趋势灵活配置混合C", "0.37", "81.37", "3.66", "0.0135", 5)
    ,@(25, "009741", "博时研究臻选三年持有期灵活配置混合C", "0.45", "80.20", "2.72", "0.0122", 8)
    ,@(26, "013798", "博时优质鑫选一年持有期混合C", "0.37", "79.82", "2.88", "0.0107", 9)
    ,@(27, "003456", "信澳新目标灵活配置混合", "0.44", "51.24", "2.03", "0.0089", 1)
    ,@(28, "960008", "景顺长城核心竞争力混合H", "0.33", "85.33", "2.56", "0.0084", 9)
    ,@(29, "005878", "博时产业新动力灵活配置混合C", "0.17", "89.28", "3.81", "0.0065", 5)
    ,@(30, "015752", "景顺长城核心招景混合C", "0.01", "89.61", "2.92", "0.0003", 9)
)

foreach ($row in $newSheetRows) {
    $r = 2 + [int]$row[0]
    $newSheet.Cells.Item($r, 1).Value = [int]$row[0]
    Set-TextCell $newSheet.Cells.Item($r, 2) $row[1]
    Set-TextCell $newSheet.Cells.Item($r, 3) $row[2]
    Set-TextCell $newSheet.Cells.Item($r, 4) $row[3]
    Set-TextCell $newSheet.Cells.Item($r, 5) $row[4]
    Set-TextCell $newSheet.Cells.Item($r, 6) $row[5]
    Set-TextCell $newSheet.Cells.Item($r, 7) $row[6]
    $newSheet.Cells.Item($r, 8).Value = [int]$row[7]
}

# The template sheet had 50 data rows (1 header + 49); our data only needs
# 32 (1 header + 31). Clear the leftover rows so the sheet's used range
# shrinks back down to A1:H32.
$newSheet.Range("A33:H50").Clear()

# Restore the original active sheet/tab (copying a sheet makes the new
# copy the active one, which would otherwise leave a stray
# tabSelected="1" on it and flip the workbook's activeTab).
$summary.Activate()

Write-Host "done"
